$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

$ws.Range("B4").Value = 2050
$ws.Range("B13").Value = 8000
$ws.Range("B19").Value = 2000

$ws.Range("B13").Select()
$ws.Activate()
